# Cycle 2 trials updated
#
# This script reproduces the "Cycle 2 trials updated" commit on the
# leaderboard worksheet:
#
#   1. The existing row 143 had the character name misspelled as
#      "YunJin" in column K; it is corrected to "Yunjin" (the spelling
#      already used elsewhere in the sheet). Once nothing references the
#      old "YunJin" text any more it naturally disappears from
#      sharedStrings.xml when the workbook is saved.
#   2. Two brand-new Cycle 2 trial rows (144 and 145) are appended below
#      the existing leaderboard table, each with a player name, cycle
#      number, video link, completion time and the 8 characters used.
#   3. The worksheet selection is left on D149, matching where the user's
#      cursor ended up after entering the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix existing row 143: "YunJin" -> "Yunjin" ----------------------
$ws.Range("K143").Value2 = "Yunjin"

# --- 2a. New row 144 -----------------------------------------------------
$ws.Range("A144").Value2 = "Not Logic"
$ws.Range("B144").Value2 = 1
$ws.Range("C144").Value2 = "https://youtu.be/tSLkS8D5cGY?si=gF0Ghf4dwNfLMSBb"
$ws.Range("D144").Value2 = 0.14930555555555555
$ws.Range("E144").Value2 = "Nilou"
$ws.Range("F144").Value2 = "Nahida"
$ws.Range("G144").Value2 = "Kokomi"
$ws.Range("H144").Value2 = "Collei"
$ws.Range("I144").Value2 = "Hutao"
$ws.Range("J144").Value2 = "Yelan"
$ws.Range("K144").Value2 = "Mona"
$ws.Range("L144").Value2 = "Zhongli"

# --- 2b. New row 145 -----------------------------------------------------
$ws.Range("A145").Value2 = "Kenny L"
$ws.Range("B145").Value2 = 1
$ws.Range("C145").Value2 = "https://youtu.be/fQ-0ucxzx7Y"
$ws.Range("D145").Value2 = 0.25694444444444448
$ws.Range("E145").Value2 = "Xingqiu"
$ws.Range("F145").Value2 = "Heizhou"
$ws.Range("G145").Value2 = "Rosaria"
$ws.Range("H145").Value2 = "Layla"
$ws.Range("I145").Value2 = "Keqing"
$ws.Range("J145").Value2 = "Fischl"
$ws.Range("K145").Value2 = "Sucrose"
$ws.Range("L145").Value2 = "Yaoyao"

# --- Carry the table's existing look down into the new rows -------------
# Column C (video link) and column D (time) both have their own
# border/number-format styling throughout the table; copy it down from
# the row directly above so rows 144/145 match rows 2-143 exactly.
try {
    $ws.Range("C143:D143").Copy()
    $ws.Range("C144:D145").PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = 0
} catch {
    # Fall back to setting the time format explicitly if PasteSpecial
    # isn't available for some reason - the important thing is the data.
    $ws.Range("D144:D145").NumberFormat = "h:mm"
}

# Make sure the D column keeps its time display even if PasteSpecial
# above did something unexpected.
$ws.Range("D144:D145").NumberFormat = "h:mm"

# --- 3. Leave the selection where the user finished typing --------------
$ws.Range("D149").Select() | Out-Null
